$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.099.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.68%  '

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.894.50'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.70%  '

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.006'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.33%  '

# Row 6: USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.22%  '

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4692'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.37%  '

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4009'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.32%  '

# Row 9: OKB
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.51'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '

# Row 10: Dogecoin
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07966'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.23%  '

# Row 11: Polygon
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9869'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '

# Row 12: Solana
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.39'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.67%  '

# Row 13: WrappedEther
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.871.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.32%  '

# Row 14: Polkadot
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.821'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '

# Row 15: Chainlink
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.012'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.56%  '

# Row 16: BinanceUSD
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.007'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.46%  '

# Row 17: Litecoin
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.73'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.45%  '

# Row 18: TRON
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06614'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.62%  '

# Row 19: ShibaInu
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001021'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.44%  '

# Row 20: Avalanche
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.27%  '

# Row 21: Dai
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '

# Row 22: WrappedBTC
$ws.Range("B22").Value = 'WrappedBTC'
$ws.Range("C22").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.153.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.87%  '

# Row 23: Uniswap
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.490'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.59%  '

# Row 24: Cosmos
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.39%  '

# Row 25: Toncoin
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.201'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.20%  '

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.123.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '

# Row 27: Monero
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.44%  '

# Row 28: EthereumClassic
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.31%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.967'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +9.19%  '

# Row 30: LidoDAOToken
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.073'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.27%  '

# Row 31: BitcoinCash
$ws.Range("B31").Value = 'BitcoinCash'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '116.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.60%  '

# Row 32: ImmutableX
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.055'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.63%  '

# Row 33: Stellar
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09431'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.23%  '

# Row 34: HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.560'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.11%  '

# Row 35: ARBITRUM
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.394'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.30%  '

# Row 36: Filecoin
$ws.Range("B36").Value = 'Filecoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.315'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.69%  '

# Row 37: Hedera
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06045'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '

# Row 38: VeChain
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02231'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.34%  '

# Row 39: TrustWalletToken
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.169'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.37%  '

# Row 40: FraxShare
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.015'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.68%  '

# Row 41: TheSandbox
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5776'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.01%  '

# Row 42: Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1816'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '

# Row 43: RenderToken
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.468'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +9.16%  '

# Row 44: Aptos
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.989'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.52%  '

# Row 45: Cronos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07666'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.03%  '

# Row 46: WEMIXToken
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.250'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.52%  '

# Row 47: EnergySwap
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.42%  '

# Row 48: Decentraland
$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5447'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.08%  '

# Row 49: NEARProtocol
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.891'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.26%  '

# Row 50: Quant
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.78'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.44%  '

# Row 51: WOONetwork
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '43.39'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.75%  '
